$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = 'Gennaro Bullo'
$ws.Range("B18").Value = 'Raffaele Prezzi  | Hellas Lazio'
$ws.Range("C18").Value = 'Andrea Conzatti | FC SAVIGNANO'
$ws.Range("D18").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("E18").Value = 'Filippo Benetti | I Magnifici'
$ws.Range("F18").Value = 'Mattia Bertolini | QUEI STRAZI'

$ws.Range("A19").Value = 'Zó'
$ws.Range("B19").Value = 'Alberto Cerisara | SHARK ATTACK'
$ws.Range("C19").Value = 'Andrea Conzatti | FC SAVIGNANO'
$ws.Range("D19").Value = 'Federico Nicolodi | U.SGUARNA'
$ws.Range("E19").Value = 'Thomas Perenzoni | CGB Gamberoni'
$ws.Range("F19").Value = 'Riccardo baldo | wanda tim'

$ws.Range("A20").Value = 'Tommibega'
$ws.Range("B20").Value = 'Samuele Kettmaier | A.C.DENTI'
$ws.Range("C20").Value = 'Andrea Menolli | SdrumALA'
$ws.Range("D20").Value = 'Luca Frasca | Clitoriders'
$ws.Range("E20").Value = 'Carlo  Stedile | Mai una gioia'
$ws.Range("F20").Value = 'Sayf Brik | A.C.DENTI'

$ws.Range("A21").Value = 'APERITIVO AL MOZART'
$ws.Range("B21").Value = 'Samuele Kettmaier | A.C.DENTI'
$ws.Range("C21").Value = 'Andrea Conzatti | FC SAVIGNANO'
$ws.Range("D21").Value = 'Andrea Menolli | SdrumALA'
$ws.Range("E21").Value = 'Giacomo  Gasparini  | Mai una gioia'
$ws.Range("F21").Value = 'Alessio Zandonai | SBARX'

$ws.Range("A22").Value = 'Biker mice'
$ws.Range("B22").Value = 'Nicolas Giordani  | FC SAVIGNANO'
$ws.Range("C22").Value = 'Federico Fasanelli | SBARX'
$ws.Range("D22").Value = 'Matteo Diener | U.SGUARNA'
$ws.Range("E22").Value = 'Filippo Benetti | I Magnifici'
$ws.Range("F22").Value = 'Alessio Debiasi | Mai una gioia'

$ws.Range("A23").Value = 'Niccoló Orsi'
$ws.Range("B23").Value = 'Lorenzo Canali | CGB Gamberoni'
$ws.Range("C23").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("D23").Value = 'Nicholas Marzadro | SBARX'
$ws.Range("E23").Value = 'Marco Sartorelli | Modium'
$ws.Range("F23").Value = 'Mattia Tezzele | U.SGUARNA'

$ws.Range("A24").Value = 'ah ma è ronco '
$ws.Range("B24").Value = 'Matteo Zanlucchi | SBARX'
$ws.Range("C24").Value = 'Carlo  Stedile | Mai una gioia'
$ws.Range("D24").Value = 'Nicolo  Speziali | FC GORILLAZ'
$ws.Range("E24").Value = 'Matteo Simoncelli | IMONTAGNA'
$ws.Range("F24").Value = 'Francesco Cristoforetti | Vigili del Fusto'

$ws.Range("A25").Value = 'Jacopo Chemini'
$ws.Range("B25").Value = 'Stefano  Tita | Clitoriders'
$ws.Range("C25").Value = 'Daniel Pedrotti | IMONTAGNA'
$ws.Range("D25").Value = 'Michele Merighi | Clitoriders'
$ws.Range("E25").Value = 'maikol  azocar | Mai una gioia'
$ws.Range("F25").Value = 'Emanuele Toss | 4SINS'

$ws.Range("A26").Value = 'Senia Lucrezia'
$ws.Range("B26").Value = 'Elia Battisti | U.SGUARNA'
$ws.Range("C26").Value = 'Daniel Pedrotti | IMONTAGNA'
$ws.Range("D26").Value = 'Michele Merighi | Clitoriders'
$ws.Range("E26").Value = 'Gabriele Verona | CGB Gamberoni'
$ws.Range("F26").Value = 'Gianni Sala | FC SALAGIARDINI'

$ws.Range("A27").Value = 'Davide Rosà'
$ws.Range("B27").Value = 'Elia Barozzi | I Magnifici'
$ws.Range("C27").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("D27").Value = 'Sebastiano Zoller | CGB Gamberoni'
$ws.Range("E27").Value = 'Filippo Benetti | I Magnifici'
$ws.Range("F27").Value = 'Mattia Tezzele | U.SGUARNA'

$ws.Range("A28").Value = 'Nazarena Raos'
$ws.Range("B28").Value = 'Alberto Cerisara | SHARK ATTACK'
$ws.Range("C28").Value = 'Andrea Gober | U.SGUARNA'
$ws.Range("D28").Value = 'Roberto Barozzi | demobusters'
$ws.Range("E28").Value = 'Andrea Conzatti | FC SAVIGNANO'
$ws.Range("F28").Value = 'Gabriele Lasta | RSA United'

$ws.Range("A29").Value = 'LORENZA SIMONCELLI'
$ws.Range("B29").Value = 'Riccardo Versini | Modium'
$ws.Range("C29").Value = 'Davide Simoncelli | Avanzi'
$ws.Range("D29").Value = 'Carlo  Stedile | Mai una gioia'
$ws.Range("E29").Value = 'Federico Fasanelli | SBARX'
$ws.Range("F29").Value = 'Emanuele Miorandi | Rita Levi’s'

$ws.Range("A30").Value = 'Thomas Cavagna'
$ws.Range("B30").Value = 'Thomas Debiasi | Mai una gioia'
$ws.Range("C30").Value = 'Thomas Cavagna | Mai una gioia'
$ws.Range("D30").Value = 'Luca Frasca | Clitoriders'
$ws.Range("E30").Value = 'Federico Nicolodi | U.SGUARNA'
$ws.Range("F30").Value = 'Davide  Bazzano | IMONTAGNA'

$ws.Range("A31").Value = 'Nicholas Marzadro'
$ws.Range("B31").Value = 'Matteo Zanlucchi | SBARX'
$ws.Range("C31").Value = 'Matteo Diener | U.SGUARNA'
$ws.Range("D31").Value = 'Federico Manica | IMONTAGNA'
$ws.Range("E31").Value = 'Filippo Benetti | I Magnifici'
$ws.Range("F31").Value = 'Alessandro Fanti | FC SALAGIARDINI'

$ws.Range("A32").Value = 'Valentina Perghem '
$ws.Range("B32").Value = 'Matteo Zanlucchi | SBARX'
$ws.Range("C32").Value = 'Luca Frasca | Clitoriders'
$ws.Range("D32").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("E32").Value = 'Michele Leonardi | Rita Levi’s'
$ws.Range("F32").Value = 'Matteo Giovannella | Bevem4tut'

$ws.Range("A33").Value = 'Emanuele Miorandi'
$ws.Range("B33").Value = 'Alberto Cerisara | SHARK ATTACK'
$ws.Range("C33").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("D33").Value = 'Filippo Benetti | I Magnifici'
$ws.Range("E33").Value = 'Alessandro Maffei | FC SAVIGNANO'
$ws.Range("F33").Value = 'Moris Benedetti | Gli Introvabili'

$ws.Range("A34").Value = 'Jasmine Scottini'
$ws.Range("B34").Value = 'Elia Battisti | U.SGUARNA'
$ws.Range("C34").Value = 'Federico Fasanelli | SBARX'
$ws.Range("D34").Value = 'Filippo Benetti | I Magnifici'
$ws.Range("E34").Value = 'Luca Frasca | Clitoriders'
$ws.Range("F34").Value = 'Francesco Cristoforetti | Vigili del Fusto'

$ws.Range("A35").Value = 'Michele Leonardi '
$ws.Range("B35").Value = 'Elia Battisti | U.SGUARNA'
$ws.Range("C35").Value = 'Filippo Benetti | I Magnifici'
$ws.Range("D35").Value = 'Matteo Diener | U.SGUARNA'
$ws.Range("E35").Value = 'Luca Tonolli | Rita Levi’s'
$ws.Range("F35").Value = 'Sayf Brik | A.C.DENTI'

$ws.Range("A36").Value = 'Mattia Spagnolli '
$ws.Range("B36").Value = 'Lorenzo Canali | CGB Gamberoni'
$ws.Range("C36").Value = 'Filippo Benetti | I Magnifici'
$ws.Range("D36").Value = 'Riccardo Zaffoni | U.SGUARNA'
$ws.Range("E36").Value = 'Sebastiano Zoller | CGB Gamberoni'
$ws.Range("F36").Value = 'Andrea Giordani | Clitoriders'

$ws.Range("A37").Value = 'Nicolas Giordani'
$ws.Range("B37").Value = 'Elia Battisti | U.SGUARNA'
$ws.Range("C37").Value = 'Filippo Benetti | I Magnifici'
$ws.Range("D37").Value = 'Federico Mortillaro | Clitoriders'
$ws.Range("E37").Value = 'Federico  Zanini | A.C.DENTI'
$ws.Range("F37").Value = 'Geremia  Carollo | FC SAVIGNANO'

$ws.Range("A38").Value = 'Michela Menghini'
$ws.Range("B38").Value = 'Nicolas Giordani  | FC SAVIGNANO'
$ws.Range("C38").Value = 'Andrea Conzatti | FC SAVIGNANO'
$ws.Range("D38").Value = 'Matteo Mazzola | GREP'
$ws.Range("E38").Value = 'Alessio Farinati | Pinguini Trentini'
$ws.Range("F38").Value = 'Emanuele  valduga | wanda tim'

$ws.Range("A39").Value = 'JR'
$ws.Range("B39").Value = 'Elia Battisti | U.SGUARNA'
$ws.Range("C39").Value = 'Luca Tonolli | Rita Levi’s'
$ws.Range("D39").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("E39").Value = 'Federico Nicolodi | U.SGUARNA'
$ws.Range("F39").Value = 'Lorenzo Mori` | Hellas Lazio'

$ws.Range("A40").Value = 'Anna Zandonati'
$ws.Range("B40").Value = 'Nicolas Giordani  | FC SAVIGNANO'
$ws.Range("C40").Value = 'Filippo Benetti | I Magnifici'
$ws.Range("D40").Value = 'Federico Fasanelli | SBARX'
$ws.Range("E40").Value = 'Luca Frasca | Clitoriders'
$ws.Range("F40").Value = 'Alessio Debiasi | Mai una gioia'

$ws.Range("A41").Value = 'Matteo Alberti'
$ws.Range("B41").Value = 'Thomas Debiasi | Mai una gioia'
$ws.Range("C41").Value = 'Nadir  chtioui | Mai una gioia'
$ws.Range("D41").Value = 'Michael Bertè  | A.C.DENTI'
$ws.Range("E41").Value = 'Andreas Galli | SdrumALA'
$ws.Range("F41").Value = 'Lorenzo Zuani | I Magnifici'

$ws.Range("A42").Value = 'Leo Parisi '
$ws.Range("B42").Value = 'Daniele Dalbosco | IMONTAGNA'
$ws.Range("C42").Value = 'Leonardo Viola | SHARK ATTACK'
$ws.Range("D42").Value = 'Geremia  Carollo | FC SAVIGNANO'
$ws.Range("E42").Value = 'Luca Frasca | Clitoriders'
$ws.Range("F42").Value = 'Moris Benedetti | Gli Introvabili'

$ws.Range("A43").Value = 'Giacomo Gasparini'
$ws.Range("B43").Value = 'Mattia Kaiserman | Gli Introvabili'
$ws.Range("C43").Value = 'Federico  Andreis | IMONTAGNA'
$ws.Range("D43").Value = 'Michele Merighi | Clitoriders'
$ws.Range("E43").Value = 'Geremia  Carollo | FC SAVIGNANO'
$ws.Range("F43").Value = 'Alessio Delli Compagni | SdrumALA'

$ws.Range("A44").Value = 'I magnifici 2.0'
$ws.Range("B44").Value = 'Stefano  Tita | Clitoriders'
$ws.Range("C44").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("D44").Value = 'Pietro  Gasparini | Mai una gioia'
$ws.Range("E44").Value = 'Geremia  Carollo | FC SAVIGNANO'
$ws.Range("F44").Value = 'Mattia Tezzele | U.SGUARNA'

$ws.Range("A45").Value = 'Raffaele Prezzi'
$ws.Range("B45").Value = 'Thomas Debiasi | Mai una gioia'
$ws.Range("C45").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("D45").Value = 'Sebastiano Zoller | CGB Gamberoni'
$ws.Range("E45").Value = 'Andreas Galli | SdrumALA'
$ws.Range("F45").Value = 'Mattia Tezzele | U.SGUARNA'

$ws.Range("A46").Value = 'Danny Giordani'
$ws.Range("B46").Value = 'Nicolas Giordani  | FC SAVIGNANO'
$ws.Range("C46").Value = 'Danny Giordani | I Magnifici'
$ws.Range("D46").Value = 'Luca Frasca | Clitoriders'
$ws.Range("E46").Value = 'Alessandro Maffei | FC SAVIGNANO'
$ws.Range("F46").Value = 'Gentian Capa | Power Ginger'

$ws.Range("A47").Value = 'Daniele Ruzzenenti'
$ws.Range("B47").Value = 'Elia Battisti | U.SGUARNA'
$ws.Range("C47").Value = 'Michele Merighi | Clitoriders'
$ws.Range("D47").Value = 'Giacomo  Gasparini  | Mai una gioia'
$ws.Range("E47").Value = 'Manuel Emanuelli | SdrumALA'
$ws.Range("F47").Value = 'Moris Benedetti | Gli Introvabili'

$ws.Range("A48").Value = 'Giovanni Simoncelli'
$ws.Range("B48").Value = 'Stefano  Tita | Clitoriders'
$ws.Range("C48").Value = 'Daniele Feltrinelli | Rita Levi’s'
$ws.Range("D48").Value = 'Sebastiano Zoller | CGB Gamberoni'
$ws.Range("E48").Value = 'Andrea  Pedrotti | IMONTAGNA'
$ws.Range("F48").Value = 'Alessandro Galvagni | Hellas Lazio'

$ws.Range("A49").Value = 'Riccardo Briosi'
$ws.Range("B49").Value = 'Riccardo Versini | Modium'
$ws.Range("C49").Value = 'Andrea Conzatti | FC SAVIGNANO'
$ws.Range("D49").Value = 'Luca Frasca | Clitoriders'
$ws.Range("E49").Value = 'Gianni Sala | FC SALAGIARDINI'
$ws.Range("F49").Value = 'Christian Torboli | 4SINS'

$ws.Range("A50").Value = 'Francesco Passuello'
$ws.Range("B50").Value = 'Gabriel Melis | demobusters'
$ws.Range("C50").Value = 'Mattia Baldessarini | SHARK ATTACK'
$ws.Range("D50").Value = 'Federico Mortillaro | Clitoriders'
$ws.Range("E50").Value = 'Mattia Tezzele | U.SGUARNA'
$ws.Range("F50").Value = 'Thomas Pontillo | Gli Introvabili'

$ws.Range("A51").Value = 'Feltri'
$ws.Range("B51").Value = 'Alberto Cerisara | SHARK ATTACK'
$ws.Range("C51").Value = 'Federico  Zanini | A.C.DENTI'
$ws.Range("D51").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("E51").Value = 'Filippo Benetti | I Magnifici'
$ws.Range("F51").Value = 'Andrea Giordani | Clitoriders'

$ws.Range("A52").Value = 'Benny'
$ws.Range("B52").Value = 'Matteo Pilati | Pinguini Trentini'
$ws.Range("C52").Value = 'Matteo Simoncelli | IMONTAGNA'
$ws.Range("D52").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("E52").Value = 'Riccardo Barbiero | Rita Levi’s'
$ws.Range("F52").Value = 'Andrea Giordani | Clitoriders'

$ws.Range("A53").Value = 'Carlotta '
$ws.Range("B53").Value = 'Daniele Dalbosco | IMONTAGNA'
$ws.Range("C53").Value = 'Andrea Bellini | Nazzzionale ferrovieri'
$ws.Range("D53").Value = 'Luca Giordani | SHARK ATTACK'
$ws.Range("E53").Value = 'Federico  Andreis | IMONTAGNA'
$ws.Range("F53").Value = 'Andrea Anzelini | GREP'

$ws.Range("A54").Value = 'Bruno 🐻'
$ws.Range("B54").Value = 'Nicolas Giordani  | FC SAVIGNANO'
$ws.Range("C54").Value = 'Marco Sala | IMONTAGNA'
$ws.Range("D54").Value = 'Danny Giordani | I Magnifici'
$ws.Range("E54").Value = 'Riccardo Zaffoni | U.SGUARNA'
$ws.Range("F54").Value = 'Alessio Debiasi | Mai una gioia'

$ws.Range("A55").Value = 'Davide Scarperi'
$ws.Range("B55").Value = 'Alessandro  Ruele  | FC GORILLAZ'
$ws.Range("C55").Value = 'Andrea Conzatti | FC SAVIGNANO'
$ws.Range("D55").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("E55").Value = 'Daniel Pedrotti | IMONTAGNA'
$ws.Range("F55").Value = 'Emanuele  valduga | wanda tim'

$ws.Range("A56").Value = 'Davide Raffaelli '
$ws.Range("B56").Value = 'Alberto Cerisara | SHARK ATTACK'
$ws.Range("C56").Value = 'Daniele Feller | GREP'
$ws.Range("D56").Value = 'Matteo Simoncelli | IMONTAGNA'
$ws.Range("E56").Value = 'Leonardo Viola | SHARK ATTACK'
$ws.Range("F56").Value = 'Alessio Debiasi | Mai una gioia'

$ws.Range("A57").Value = 'Sebastiano Zoller'
$ws.Range("B57").Value = 'Lorenzo Canali | CGB Gamberoni'
$ws.Range("C57").Value = 'Sebastiano Zoller | CGB Gamberoni'
$ws.Range("D57").Value = 'Andrea  Roveda  | Pinguini Trentini'
$ws.Range("E57").Value = 'Michele Merighi | Clitoriders'
$ws.Range("F57").Value = 'Alessio Debiasi | Mai una gioia'

$ws.Range("A58").Value = 'Gabriele Gottardi'
$ws.Range("B58").Value = 'Federico Zoller | GREP'
$ws.Range("C58").Value = 'Andrea  Roveda  | Pinguini Trentini'
$ws.Range("D58").Value = 'Luca Perenzoni | CGB Gamberoni'
$ws.Range("E58").Value = 'Michele Merighi | Clitoriders'
$ws.Range("F58").Value = 'Alessio  Giordano  | FC Schalke 104'

$ws.Range("A59").Value = 'Riccardo Zeni'
$ws.Range("B59").Value = 'Elia Barozzi | I Magnifici'
$ws.Range("C59").Value = 'Sebastiano Zoller | CGB Gamberoni'
$ws.Range("D59").Value = 'Leonardo Viola | SHARK ATTACK'
$ws.Range("E59").Value = 'Andrea  Roveda  | Pinguini Trentini'
$ws.Range("F59").Value = 'Davide  Bazzano | IMONTAGNA'

$ws.Range("A60").Value = 'Davide Zeni'
$ws.Range("B60").Value = 'Lorenzo Canali | CGB Gamberoni'
$ws.Range("C60").Value = 'Andrea  Roveda  | Pinguini Trentini'
$ws.Range("D60").Value = 'Alessio Bragagna | SHARK ATTACK'
$ws.Range("E60").Value = 'Luca Perenzoni | CGB Gamberoni'
$ws.Range("F60").Value = 'Jacopo  Chemini | IMONTAGNA'

$ws.Range("A61").Value = 'Riccardo Barbiero'
$ws.Range("B61").Value = 'Elia Battisti | U.SGUARNA'
$ws.Range("C61").Value = 'Riccardo Barbiero | Rita Levi’s'
$ws.Range("D61").Value = 'Leonardo Viola | SHARK ATTACK'
$ws.Range("E61").Value = 'Marco Sala | IMONTAGNA'
$ws.Range("F61").Value = 'Moris Benedetti | Gli Introvabili'

$ws.Range("A62:F67").Clear()
